# Resonator Syn - Initial Sweep (optimal Bridge=0.28)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Resonator Synthesis")

# --- ATTEMPT 2 banner + mini summary table (C76:D78) ---
$ws.Range("C76").Value = " ======================= ATTEMPT 2 ========================"

$ws.Range("C77").Value = "Bridge"
$ws.Range("D77").Value = "Flatness"

$ws.Range("C78").Value = 0.2
$ws.Range("D78").Value = "N/A"
$ws.Range("C79").Value = 0.22
$ws.Range("C80").Value = 0.24
$ws.Range("C81").Value = 0.26
$ws.Range("C82").Value = 0.28
$ws.Range("C83").Value = 0.3

# --- Sweep table #1 (J77:L83) ---
$ws.Range("J77").Value = "Bridge"
$ws.Range("K77").Value = "Leg"
$ws.Range("L77").Value = "f0"

$ws.Range("J78").Value = 0.22
$ws.Range("K78").Value = 0.3
$ws.Range("L78").Value = 4.0859049207672999

$ws.Range("J79").Value = 0.22
$ws.Range("K79").Value = 0.32
$ws.Range("L79").Value = 3.8590492076730998

$ws.Range("J80").Value = 0.22
$ws.Range("K80").Value = 0.34
$ws.Range("L80").Value = 3.8590492076730998

$ws.Range("J81").Value = 0.22
$ws.Range("K81").Value = 0.36
$ws.Range("L81").Value = 3.8990825688072999

$ws.Range("J82").Value = 0.22
$ws.Range("K82").Value = 0.38
$ws.Range("L82").Value = 3.7589658048374002

$ws.Range("J83").Value = 0.22
$ws.Range("K83").Value = 0.4
$ws.Range("L83").Value = 3.3619683069224

# --- Sweep table #2 (N77:P83) ---
$ws.Range("N77").Value = "Bridge"
$ws.Range("O77").Value = "Leg"
$ws.Range("P77").Value = "f0"

$ws.Range("N78").Value = 0.24
$ws.Range("O78").Value = 0.3
$ws.Range("P78").Value = 4.0859049207672999

$ws.Range("N79").Value = 0.24
$ws.Range("O79").Value = 0.32
$ws.Range("P79").Value = 3.8590492076730998

$ws.Range("N80").Value = 0.24
$ws.Range("O80").Value = 0.34
$ws.Range("P80").Value = 3.8590492076730998

$ws.Range("N81").Value = 0.24
$ws.Range("O81").Value = 0.36
$ws.Range("P81").Value = 3.8990825688072999

$ws.Range("N82").Value = 0.24
$ws.Range("O82").Value = 0.38
$ws.Range("P82").Value = 3.7589658048374002

$ws.Range("N83").Value = 0.24
$ws.Range("O83").Value = 0.4
$ws.Range("P83").Value = 3.3619683069224

# --- Make "Resonator Synthesis" the active sheet/view, matching the saved cursor ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 70
$ws.Range("N84").Select()
